$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.192.57"
$ws.Range("E2").Value = "'  -2.85%  "
$ws.Range("D3").Value = "'1.928.88"
$ws.Range("E3").Value = "'  -1.56%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'323.15"
$ws.Range("E5").Value = "'  -1.07%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "'  +0.34%  "
$ws.Range("D7").Value = "'0.4741"
$ws.Range("E7").Value = "'  -4.28%  "
$ws.Range("D8").Value = "'0.4052"
$ws.Range("E8").Value = "'  -3.34%  "
$ws.Range("D9").Value = "'53.62"
$ws.Range("E9").Value = "'  +1.61%  "
$ws.Range("D10").Value = "'0.08507"
$ws.Range("E10").Value = "'  -7.41%  "
$ws.Range("D11").Value = "'1.050"
$ws.Range("E11").Value = "'  -4.07%  "
$ws.Range("E12").Value = "'  -2.51%  "
$ws.Range("D13").Value = "'1.970.83"
$ws.Range("E13").Value = "'  -2.76%  "
$ws.Range("D14").Value = "'7.522"
$ws.Range("E14").Value = "'  -3.81%  "
$ws.Range("D15").Value = "'6.128"
$ws.Range("E15").Value = "'  -4.71%  "
$ws.Range("E16").Value = "'  +0.46%  "
$ws.Range("D17").Value = "'89.91"
$ws.Range("E17").Value = "'  -1.43%  "
$ws.Range("D18").Value = "'0.00001068"
$ws.Range("E18").Value = "'  -2.47%  "
$ws.Range("D19").Value = "'0.06592"
$ws.Range("E19").Value = "'  -1.49%  "
$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = "'  -5.43%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("E21").Value = "'  +0.48%  "
$ws.Range("D22").Value = "'5.777"
$ws.Range("E22").Value = "'  -2.81%  "
$ws.Range("D23").Value = "'28.241.60"
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "'  -4.68%  "
$ws.Range("D25").Value = "'2.295"
$ws.Range("E25").Value = "'  +1.38%  "
$ws.Range("D26").Value = "'2.203.36"
$ws.Range("E26").Value = "'  -2.20%  "
$ws.Range("D27").Value = "'154.90"
$ws.Range("E27").Value = "'  -0.41%  "
$ws.Range("D28").Value = "'20.16"
$ws.Range("E28").Value = "'  -1.88%  "
$ws.Range("D29").Value = "'2.170"
$ws.Range("E29").Value = "'  -3.25%  "
$ws.Range("D30").Value = "'5.769"
$ws.Range("E30").Value = "'  -7.84%  "
$ws.Range("D31").Value = "'123.88"
$ws.Range("E31").Value = "'  -1.74%  "
$ws.Range("D32").Value = "'0.9818"
$ws.Range("E32").Value = "'  -5.53%  "
$ws.Range("D33").Value = "'0.09604"
$ws.Range("E33").Value = "'  -2.07%  "
$ws.Range("D34").Value = "'1.445"
$ws.Range("E34").Value = "'  -4.87%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.587"
$ws.Range("E35").Value = "'  -3.94%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'3.654"
$ws.Range("E36").Value = "'  -0.71%  "
$ws.Range("D37").Value = "'9.279"
$ws.Range("E37").Value = "'  +3.46%  "
$ws.Range("D38").Value = "'0.02324"
$ws.Range("E38").Value = "'  -3.86%  "
$ws.Range("D39").Value = "'0.06174"
$ws.Range("E39").Value = "'  -2.61%  "
$ws.Range("E40").Value = "'  -6.65%  "
$ws.Range("D41").Value = "'0.6183"
$ws.Range("E41").Value = "'  -3.68%  "
$ws.Range("D42").Value = "'11.11"
$ws.Range("E42").Value = "'  -2.17%  "
$ws.Range("D43").Value = "'1.008"
$ws.Range("E43").Value = "'  +0.35%  "
$ws.Range("D44").Value = "'0.1910"
$ws.Range("E44").Value = "'  -3.17%  "
$ws.Range("D45").Value = "'1.320"
$ws.Range("E45").Value = "'  -2.99%  "
$ws.Range("D46").Value = "'0.5894"
$ws.Range("E46").Value = "'  -4.94%  "
$ws.Range("D47").Value = "'12.84"
$ws.Range("E47").Value = "'  -3.57%  "
$ws.Range("D48").Value = "'2.042"
$ws.Range("E48").Value = "'  -6.60%  "
$ws.Range("D49").Value = "'3.405"
$ws.Range("E49").Value = "'  -1.56%  "
$ws.Range("D50").Value = "'0.06784"
$ws.Range("E50").Value = "'  -3.80%  "
$ws.Range("D51").Value = "'1.088"
$ws.Range("E51").Value = "'  -2.18%  "
